# Finalized definition of virtual environment
# Updates the calibration outputs across the AR / SETAR / GARCH / TARCH /
# AR_TARCH sheets (new coefficients + refreshed abs_epsi_autocorr series).

$wb = $excel.ActiveWorkbook

# --- AR sheet ---------------------------------------------------------
$wsAR = $wb.Worksheets.Item("AR")
$wsAR.Range("B2").Value = -0.02935319036800062
$wsAR.Range("B3").Value = 0.8148737245007158
$wsAR.Range("B4").Value = 3.716600370622494
$wsAR.Range("B5").Value = "[1.0, 0.2234867022365352, 0.2172810577187567, 0.22470527548345198, 0.1807468192024304, 0.43315426219758657, 0.15416538919348818, 0.1918917660781605, 0.15586654148211257, 0.15558596961612492, 0.17621703721857068, 0.14127464040343926, 0.14492838670492644, 0.1541247001400519, 0.20123628523550954, 0.13041877917144956, 0.1383669024519305, 0.13308786916572388, 0.1490422300069606, 0.1417354568920585]"

# --- SETAR sheet --------------------------------------------------------
$wsSETAR = $wb.Worksheets.Item("SETAR")
$wsSETAR.Range("B3").Value = 0.9504950495049505
$wsSETAR.Range("B4").Value = -0.5027905522458885
$wsSETAR.Range("B5").Value = 0.7890286697769049
$wsSETAR.Range("B6").Value = 2.815528618716
$wsSETAR.Range("B7").Value = 0.4880421130553186
$wsSETAR.Range("B8").Value = 0.7499760826518644
$wsSETAR.Range("B9").Value = 2.318611603655023
$wsSETAR.Range("B10").Value = "[1.0, 0.27058448663998214, 0.23717296025651552, 0.21254927782762872, 0.23408012621294352, 0.3183038808926737, 0.17249585387321204, 0.1458458622421585, 0.11797377222555902, 0.11979842423601327, 0.13807815108421057, 0.1237746162416953, 0.12624557104638925, 0.11011414679294489, 0.16721991662563385, 0.10719416027306998, 0.11261688152933319, 0.13592968102345418, 0.13520540254759933, 0.13016966520985004]"

# --- GARCH sheet ---------------------------------------------------------
$wsGARCH = $wb.Worksheets.Item("GARCH")
$wsGARCH.Range("B2").Value = -0.008228157979410178
$wsGARCH.Range("B3").Value = 0.02440120630422643
$wsGARCH.Range("B4").Value = 0.0558718297427872
$wsGARCH.Range("B5").Value = 0.9379161340050154
$wsGARCH.Range("B6").Value = "[1.0, 0.009800098672221034, -0.01113181124804023, -0.01559270932403285, -0.04703880614424163, 0.3610349224913171, 0.007585071214567595, 0.01710742618962009, -0.014087888964287788, -0.044030782202712365, -0.042260052597155415, -0.04526943186985814, -0.02412474422429682, -0.008491266206800558, 0.015773343561695095, -0.02009750611840592, -0.04669707803652939, -0.017175291806260533, -0.023442997254830046, -0.010455670554982286]"

# --- TARCH sheet ---------------------------------------------------------
$wsTARCH = $wb.Worksheets.Item("TARCH")
$wsTARCH.Range("B2").Value = 0.0001541063721840237
$wsTARCH.Range("B3").Value = 0.02517557180930631
$wsTARCH.Range("B4").Value = 0.06406123341689407
$wsTARCH.Range("B5").Value = 0.936541274588946
$wsTARCH.Range("B6").Value = "[1.0, 0.009534978884775292, -0.01171852898527423, -0.015695598660415283, -0.048254005867216884, 0.36063381513269327, 0.007118110659867282, 0.015843386456236997, -0.014073954113537082, -0.04478499240823882, -0.04226127077380087, -0.045901432876203206, -0.02419036096886461, -0.008574949564069703, 0.015756502324402375, -0.02037624246782625, -0.04760399244505802, -0.017888415808653125, -0.023994971953874393, -0.01047201053134273]"
$wsTARCH.Range("B7").Value = -0.01354160067825516

# --- AR_TARCH sheet -------------------------------------------------------
$wsARTARCH = $wb.Worksheets.Item("AR_TARCH")
$wsARTARCH.Range("B2").Value = -0.02070492103214357
$wsARTARCH.Range("B3").Value = 0.02270486645974969
$wsARTARCH.Range("B4").Value = 0.07090128650642774
$wsARTARCH.Range("B5").Value = 0.9364212174099417
$wsARTARCH.Range("B6").Value = "[0.9999999999999999, 0.04252618163271713, 0.009361987857899205, 0.03554972031141428, -0.018432362910407424, 0.2633848992421267, -0.0038824135871978176, 0.009216933163390558, -0.025987661653650622, -0.0452950779510569, -0.03622721861242116, -0.04628307394029194, -0.03569159338141229, -0.022218979549180496, 0.014290349540840715, -0.018584960992281463, -0.04466796030942388, -0.02568730662106427, -0.01751024810194324, -0.016425809126267565]"
$wsARTARCH.Range("B7").Value = -0.02662103703004695
$wsARTARCH.Range("B9").Value = 0.7911148618580267
